{"js": "// The author appended a new sentence to the end of the \"Proposal: ...\"\n// paragraph in the CQL 2 Use Cases document. The new text reads:\n//   \"  The schema also contains the AttributeValue type, which is a choice\n//   of string, date, integer, long, and boolean typed elements.\"\n// (note the two leading spaces, matching the original diff's runs).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"Proposal: The new CQLAttribute.xsd schema\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(marker) === 0) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the target paragraph starting with: \" + marker);\n}\n\nconst addition =\n  \"  The schema also contains the AttributeValue type, which is a choice \" +\n  \"of string, date, integer, long, and boolean typed elements.\";\n\n// Insert as a new run at the end of the existing paragraph (no new\n// paragraph mark is introduced \u2014 the sentence is appended inline).\ntarget.insertText(addition, Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# The author appended a new sentence to the end of the \"Proposal: ...\"\n# paragraph in the CQL 2 Use Cases document. The new text reads:\n#   \"  The schema also contains the AttributeValue type, which is a choice\n#   of string, date, integer, long, and boolean typed elements.\"\n# (note the two leading spaces).\n\n$d = $word.ActiveDocument\n\n$marker = \"Proposal: The new CQLAttribute.xsd schema\"\n$addition = \"  The schema also contains the AttributeValue type, which is a choice of string, date, integer, long, and boolean typed elements.\"\n\n$target = $null\nforeach ($para in $d.Paragraphs) {\n    if ($para.Range.Text.StartsWith($marker)) {\n        $target = $para\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the target paragraph starting with: $marker\"\n}\n\n$target.Range.InsertAfter($addition)\n"}
